$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "'3"
$ws.Range("B4").Value = "'76788"
$ws.Range("C4").Value = "'213123"
$ws.Range("D4").Value = "'2"
$ws.Range("E4").Value = "info                          "
$ws.Range("F4").Value = "'5"

$ws.Range("A5").Value = "'4"
$ws.Range("B5").Value = "'121212"
$ws.Range("C5").Value = "'56489299"
$ws.Range("D5").Value = "'5"
$ws.Range("E5").Value = "asddawsxx                     "
$ws.Range("F5").Value = "'8"

$ws.Range("A4:F5").ClearFormats()
